$d = $word.ActiveDocument

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.InsertAfter("Change ")
